# Marsframe work Sprint 2
# Update the SignIn sheet's test data (Url / Username / Password) and make
# SignIn the active sheet/tab (previously ManageListings was active).

$wb = $excel.ActiveWorkbook

$signIn = $wb.Worksheets.Item("SignIn")

$signIn.Range("A2").Value = "http://localhost:5000/"
$signIn.Range("B2").Value = "aswini.sanal@outlook.com"
$signIn.Range("C2").Value = "Testing0123*"

$signIn.Activate()
